$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 26.92947733333333
$ws.Range("H2").Value = 80.788432
$ws.Range("I2").Value = 0.9279949792877585
$ws.Range("J2").Value = 0.9279949792877585
$ws.Range("M2").Value = 0.5843273333333333
$ws.Range("N2").Value = 1.752982
$ws.Range("O2").Value = 0.007060666168730862
$ws.Range("P2").Value = 0.007060666168730861
$ws.Range("Q2").Value = 15.73562967824711
$ws.Range("R2").Value = 141.620667104224
$ws.Range("S2").Value = 0.006552262755009174
$ws.Range("T2").Value = 0.006552262755009173

# Row 3
$ws.Range("G3").Value = 26.92947733333333
$ws.Range("H3").Value = 80.788432
$ws.Range("I3").Value = 0.9279949792877585
$ws.Range("J3").Value = 0.9279949792877585
$ws.Range("O3").Value = 0.9042969637611353
$ws.Range("P3").Value = 0.9042969637611354
$ws.Range("Q3").Value = 2015.345549677252
$ws.Range("R3").Value = 18138.10994709526
$ws.Range("S3").Value = 0.8391830421554977
$ws.Range("T3").Value = 0.8391830421554978

# Row 4
$ws.Range("G4").Value = 26.92947733333333
$ws.Range("H4").Value = 80.788432
$ws.Range("I4").Value = 0.9279949792877585
$ws.Range("J4").Value = 0.9279949792877585
$ws.Range("M4").Value = 7.335874333333333
$ws.Range("O4").Value = 0.08864237007013374
$ws.Range("P4").Value = 0.08864237007013374
$ws.Range("Q4").Value = 197.5512615796818
$ws.Range("R4").Value = 1777.961354217136
$ws.Range("S4").Value = 0.0822596743772516
$ws.Range("T4").Value = 0.0822596743772516

# Row 5
$ws.Range("I5").Value = 0.04103565698374688
$ws.Range("J5").Value = 0.04103565698374688
$ws.Range("M5").Value = 0.5843273333333333
$ws.Range("N5").Value = 1.752982
$ws.Range("O5").Value = 0.007060666168730862
$ws.Range("P5").Value = 0.007060666168730861
$ws.Range("Q5").Value = 0.6958247795644444
$ws.Range("R5").Value = 6.262423016079999
$ws.Range("S5").Value = 0.0002897390749767859
$ws.Range("T5").Value = 0.0002897390749767859

# Row 6
$ws.Range("I6").Value = 0.04103565698374688
$ws.Range("J6").Value = 0.04103565698374688
$ws.Range("O6").Value = 0.9042969637611353
$ws.Range("P6").Value = 0.9042969637611354
$ws.Range("S6").Value = 0.03710842001634573
$ws.Range("T6").Value = 0.03710842001634573

# Row 7
$ws.Range("I7").Value = 0.04103565698374688
$ws.Range("J7").Value = 0.04103565698374688
$ws.Range("M7").Value = 7.335874333333333
$ws.Range("O7").Value = 0.08864237007013374
$ws.Range("P7").Value = 0.08864237007013374
$ws.Range("Q7").Value = 8.73565696779111
$ws.Range("R7").Value = 78.62091271012
$ws.Range("S7").Value = 0.003637497892424359
$ws.Range("T7").Value = 0.003637497892424359

# Row 8
$ws.Range("G8").Value = 0.8986996666666666
$ws.Range("I8").Value = 0.03096936372849452
$ws.Range("J8").Value = 0.03096936372849452
$ws.Range("M8").Value = 0.5843273333333333
$ws.Range("N8").Value = 1.752982
$ws.Range("O8").Value = 0.007060666168730862
$ws.Range("P8").Value = 0.007060666168730861
$ws.Range("Q8").Value = 0.5251347796908888
$ws.Range("R8").Value = 4.726213017217999
$ws.Range("S8").Value = 0.0002186643387449019
$ws.Range("T8").Value = 0.0002186643387449019

# Row 9
$ws.Range("G9").Value = 0.8986996666666666
$ws.Range("I9").Value = 0.03096936372849452
$ws.Range("J9").Value = 0.03096936372849452
$ws.Range("O9").Value = 0.9042969637611353
$ws.Range("P9").Value = 0.9042969637611354
$ws.Range("Q9").Value = 67.25679638316645
$ws.Range("R9").Value = 605.3111674484981
$ws.Range("S9").Value = 0.02800550158929183
$ws.Range("T9").Value = 0.02800550158929184

# Row 10
$ws.Range("G10").Value = 0.8986996666666666
$ws.Range("I10").Value = 0.03096936372849452
$ws.Range("J10").Value = 0.03096936372849452
$ws.Range("M10").Value = 7.335874333333333
$ws.Range("O10").Value = 0.08864237007013374
$ws.Range("P10").Value = 0.08864237007013374
$ws.Range("Q10").Value = 6.592747818075221
$ws.Range("R10").Value = 59.33473036267699
$ws.Range("S10").Value = 0.002745197800457788
$ws.Range("T10").Value = 0.002745197800457789
